$wb = $excel.ActiveWorkbook

# --- DatosCuenta (sheet 1) ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokePreCuatro"
$wsCuenta.Range("B2").Value = "SmokePreNameCuatro"
$wsCuenta.Range("C2").Value = 27100111
$wsCuenta.Range("D2").Value = 113

# --- DatosHogar (sheet 2) ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 633

# --- DatosMotor (sheet 3) ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMA014"
$wsMotor.Range("B2").Value = "ABC12SSMA014"
$wsMotor.Range("C2").Value = "ZAZ123SSMA014"

# --- DatosAP (sheet 4) ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200113

# Selection / active-sheet bookkeeping: the active tab moves from DatosAP (C15)
# to DatosCuenta (D2), and DatosAP's own lingering selection moves to E4.
$wsAP.Range("E4").Select() | Out-Null
$wsCuenta.Range("D2").Select() | Out-Null
